$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November")

$ws.Range("B2").Value = 1496
$ws.Range("C2").Value = 1082
$ws.Range("D2").Value = 414
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.38 : 1"

$ws.Range("B3").Value = 501
$ws.Range("C3").Value = 353
$ws.Range("D3").Value = 148
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.42 : 1"

$ws.Range("B4").Value = 1088
$ws.Range("C4").Value = 1181
$ws.Range("D4").Value = -93
$ws.Range("F4").Value = "We lent more than we borrowed"
$ws.Range("G4").Value = "0.92 : 1"

$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 107
$ws.Range("D5").Value = -67
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.37 : 1"

$ws.Range("B6").Value = 899
$ws.Range("C6").Value = 1318
$ws.Range("D6").Value = -419
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.68 : 1"

$ws.Range("B7").Value = 128
$ws.Range("C7").Value = 161
$ws.Range("D7").Value = -33
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.80 : 1"

$ws.Range("B8").Value = 133
$ws.Range("C8").Value = 185
$ws.Range("D8").Value = -52
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.72 : 1"

$ws.Range("B9").Value = 29
$ws.Range("C9").Value = 59
$ws.Range("D9").Value = -30
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.49 : 1"

$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 41
$ws.Range("D10").Value = -36
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.12 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 18
$ws.Range("D12").Value = -7
$ws.Range("F12").Value = "We lent more than we borrowed"
$ws.Range("G12").Value = "0.61 : 1"

$ws.Range("B13").Value = 75
$ws.Range("C13").Value = 41
$ws.Range("D13").Value = 34
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.83 : 1"

$ws.Range("B14").Value = 200
$ws.Range("C14").Value = 203
$ws.Range("D14").Value = -3
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.99 : 1"

$ws.Range("B15").Value = 48
$ws.Range("C15").Value = 118
$ws.Range("D15").Value = -70
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.41 : 1"

$ws.Range("B16").Value = 37
$ws.Range("C16").Value = 139
$ws.Range("D16").Value = -102
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.27 : 1"

$ws.Range("B17").Value = 552
$ws.Range("C17").Value = 350
$ws.Range("D17").Value = 202
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.58 : 1"

$ws.Range("B18").Value = 87
$ws.Range("C18").Value = 91
$ws.Range("D18").Value = -4
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.96 : 1"

$ws.Range("B19").Value = 444
$ws.Range("C19").Value = 313
$ws.Range("D19").Value = 131
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.42 : 1"

$ws.Range("B20").Value = 37
$ws.Range("C20").Value = 49
$ws.Range("D20").Value = -12
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.76 : 1"

$ws.Range("B21").Value = 500
$ws.Range("C21").Value = 271
$ws.Range("D21").Value = 229
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.85 : 1"

$ws.Range("B22").Value = 33
$ws.Range("C22").Value = 139
$ws.Range("D22").Value = -106
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.24 : 1"

$ws.Range("B23").Value = 475
$ws.Range("C23").Value = 238
$ws.Range("D23").Value = 237
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "2.00 : 1"

$ws.Range("B24").Value = 1337
$ws.Range("C24").Value = 981
$ws.Range("D24").Value = 356
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.36 : 1"

$ws.Range("B25").Value = 166
$ws.Range("C25").Value = 473
$ws.Range("D25").Value = -307
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.35 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 214
$ws.Range("C27").Value = 167
$ws.Range("D27").Value = 47
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.28 : 1"

$ws.Range("B28").Value = 112
$ws.Range("C28").Value = 81
$ws.Range("D28").Value = 31
$ws.Range("E28").Value = "We borrowerd more than we lent"
$ws.Range("G28").Value = "1.38 : 1"

$ws.Range("B29").Value = 661
$ws.Range("C29").Value = 467
$ws.Range("D29").Value = 194
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.42 : 1"

$ws.Range("B30").Value = 13
$ws.Range("C30").Value = 38
$ws.Range("D30").Value = -25
$ws.Range("F30").Value = "We lent more than we borrowed"
$ws.Range("G30").Value = "0.34 : 1"

$ws.Range("B31").Value = 63
$ws.Range("C31").Value = 238
$ws.Range("D31").Value = -175
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.26 : 1"

$ws.Range("B32").Value = 437
$ws.Range("C32").Value = 488
$ws.Range("D32").Value = -51
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.90 : 1"

$ws.Range("B33").Value = 254
$ws.Range("C33").Value = 454
$ws.Range("D33").Value = -200
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.56 : 1"

$ws.Range("B34").Value = 148
$ws.Range("C34").Value = 149
$ws.Range("D34").Value = -1
$ws.Range("F34").Value = "We lent more than we borrowed"
$ws.Range("G34").Value = "0.99 : 1"

$ws.Range("B35").Value = 855
$ws.Range("C35").Value = 1031
$ws.Range("D35").Value = -176
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.83 : 1"

$ws.Range("B36").Value = 133
$ws.Range("C36").Value = 511
$ws.Range("D36").Value = -378
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.26 : 1"

$ws.Range("B37").Value = 462
$ws.Range("C37").Value = 299
$ws.Range("D37").Value = 163
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.55 : 1"

$ws.Range("B38").Value = 31
$ws.Range("C38").Value = 135
$ws.Range("D38").Value = -104
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.23 : 1"

$ws.Range("B39").Value = 18
$ws.Range("C39").Value = 97
$ws.Range("D39").Value = -79
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.19 : 1"

$ws.Range("B40").Value = 109
$ws.Range("C40").Value = 122
$ws.Range("D40").Value = -13
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.89 : 1"

$ws.Range("B41").Value = 16
$ws.Range("C41").Value = 32
$ws.Range("D41").Value = -16
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.50 : 1"

$ws.Range("B42").Value = 11
$ws.Range("C42").Value = 29
$ws.Range("D42").Value = -18
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.38 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 84
$ws.Range("C44").Value = 86
$ws.Range("D44").Value = -2
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.98 : 1"

$ws.Range("B45").Value = 52
$ws.Range("C45").Value = 193
$ws.Range("D45").Value = -141
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.27 : 1"

$ws.Range("B46").Value = 417
$ws.Range("C46").Value = 541
$ws.Range("D46").Value = -124
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.77 : 1"

$ws.Range("B47").Value = 857
$ws.Range("C47").Value = 555
$ws.Range("D47").Value = 302
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.54 : 1"

$ws.Range("B48").Value = 172
$ws.Range("C48").Value = 578
$ws.Range("D48").Value = -406
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.30 : 1"

$ws.Range("B49").Value = 660
$ws.Range("C49").Value = 252
$ws.Range("D49").Value = 408
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "2.62 : 1"

$ws.Range("B50").Value = 690
$ws.Range("C50").Value = 524
$ws.Range("D50").Value = 166
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.32 : 1"

$ws.Range("B51").Value = 191
$ws.Range("C51").Value = 161
$ws.Range("D51").Value = 30
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.19 : 1"

$ws.Range("B52").Value = 312
$ws.Range("C52").Value = 322
$ws.Range("D52").Value = -10
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.97 : 1"

$ws.Range("B53").Value = 233
$ws.Range("C53").Value = 177
$ws.Range("D53").Value = 56
$ws.Range("E53").Value = "We borrowerd more than we lent"
$ws.Range("G53").Value = "1.32 : 1"

$ws.Range("B54").Value = 38
$ws.Range("C54").Value = 231
$ws.Range("D54").Value = -193
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.16 : 1"

$ws.Range("B55").Value = 422
$ws.Range("C55").Value = 117
$ws.Range("D55").Value = 305
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "3.61 : 1"

# Restore selection: user last selected November!E11, then returned to Yearly totals tab
$ws.Range("E11").Select()
$wb.Worksheets.Item("Yearly totals").Activate()

Write-Host "November sheet updated"